$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.19533527696793
$ws.Range("C2").Value = 0.5568513119533528
$ws.Range("J2").Value = 0.02040816326530612
$ws.Range("P2").Value = 0.1457725947521866
$ws.Range("S2").Value = 0.08163265306122448
$ws.Range("B3").Value = 0.004672897196261682
$ws.Range("C3").Value = 0.03738317757009346
$ws.Range("J3").Value = 0.04205607476635514
$ws.Range("P3").Value = 0.7570093457943925
$ws.Range("S3").Value = 0.1588785046728972
$ws.Range("J4").Value = 0.0847457627118644
$ws.Range("O4").Value = 0.01694915254237288
$ws.Range("P4").Value = 0.5423728813559322
$ws.Range("S4").Value = 0.3559322033898305
$ws.Range("B6").Value = 0.07608695652173914
$ws.Range("D6").Value = 0.01449275362318841
$ws.Range("F6").Value = 0.04710144927536232
$ws.Range("J6").Value = 0.2065217391304348
$ws.Range("O6").Value = 0.02898550724637681
$ws.Range("Q6").Value = 0.1884057971014493
$ws.Range("R6").Value = 0.07608695652173914
$ws.Range("S6").Value = 0.3623188405797101
$ws.Range("B7").Value = 0.1405405405405405
$ws.Range("D7").Value = 0.01081081081081081
$ws.Range("F7").Value = 0.04864864864864865
$ws.Range("J7").Value = 0.1135135135135135
$ws.Range("O7").Value = 0.02162162162162162
$ws.Range("Q7").Value = 0.1891891891891892
$ws.Range("R7").Value = 0.04864864864864865
$ws.Range("S7").Value = 0.4270270270270271
$ws.Range("B8").Value = 0.07205240174672489
$ws.Range("D8").Value = 0.01528384279475982
$ws.Range("E8").Value = 0.002183406113537118
$ws.Range("F8").Value = 0.07423580786026202
$ws.Range("J8").Value = 0.1157205240174673
$ws.Range("O8").Value = 0.02183406113537118
$ws.Range("Q8").Value = 0.2008733624454148
$ws.Range("R8").Value = 0.07205240174672489
$ws.Range("S8").Value = 0.425764192139738
$ws.Range("B9").Value = 0.1219512195121951
$ws.Range("D9").Value = 0.02439024390243903
$ws.Range("F9").Value = 0.08130081300813008
$ws.Range("J9").Value = 0.0975609756097561
$ws.Range("O9").Value = 0.01626016260162602
$ws.Range("Q9").Value = 0.2154471544715447
$ws.Range("R9").Value = 0.08130081300813008
$ws.Range("S9").Value = 0.3617886178861789
$ws.Range("B10").Value = 0.1163793103448276
$ws.Range("D10").Value = 0.02658045977011494
$ws.Range("E10").Value = 0.001436781609195402
$ws.Range("F10").Value = 0.08764367816091954
$ws.Range("J10").Value = 0.1063218390804598
$ws.Range("O10").Value = 0.02514367816091954
$ws.Range("Q10").Value = 0.2119252873563219
$ws.Range("R10").Value = 0.07471264367816093
$ws.Range("S10").Value = 0.3498563218390804
$ws.Range("G11").Value = 0.1192982456140351
$ws.Range("J11").Value = 0.1192982456140351
$ws.Range("K11").Value = 0.1543859649122807
$ws.Range("L11").Value = 0.5859649122807018
$ws.Range("S11").Value = 0.02105263157894737
$ws.Range("G12").Value = 0.7444444444444445
$ws.Range("J12").Value = 0.15
$ws.Range("K12").Value = 0.01111111111111111
$ws.Range("L12").Value = 0.06111111111111111
$ws.Range("S12").Value = 0.03333333333333333
$ws.Range("G13").Value = 0.5476190476190477
$ws.Range("J13").Value = 0.3571428571428572
$ws.Range("S13").Value = 0.09523809523809523
$ws.Range("F15").Value = 0.01374570446735395
$ws.Range("H15").Value = 0.1718213058419244
$ws.Range("I15").Value = 0.06872852233676977
$ws.Range("J15").Value = 0.3608247422680412
$ws.Range("K15").Value = 0.04810996563573883
$ws.Range("M15").Value = 0.01374570446735395
$ws.Range("O15").Value = 0.08934707903780069
$ws.Range("S15").Value = 0.2336769759450172
$ws.Range("F16").Value = 0.02690582959641256
$ws.Range("H16").Value = 0.1300448430493273
$ws.Range("I16").Value = 0.1076233183856502
$ws.Range("J16").Value = 0.4170403587443946
$ws.Range("K16").Value = 0.1210762331838565
$ws.Range("M16").Value = 0.02690582959641256
$ws.Range("O16").Value = 0.04484304932735426
$ws.Range("S16").Value = 0.1255605381165919
$ws.Range("F17").Value = 0.02238805970149254
$ws.Range("H17").Value = 0.1828358208955224
$ws.Range("I17").Value = 0.1156716417910448
$ws.Range("J17").Value = 0.4421641791044776
$ws.Range("K17").Value = 0.07462686567164178
$ws.Range("M17").Value = 0.01865671641791045
$ws.Range("O17").Value = 0.06529850746268656
$ws.Range("S17").Value = 0.07835820895522388
$ws.Range("F18").Value = 0.01570680628272251
$ws.Range("H18").Value = 0.1780104712041885
$ws.Range("I18").Value = 0.1256544502617801
$ws.Range("J18").Value = 0.4397905759162304
$ws.Range("K18").Value = 0.08900523560209424
$ws.Range("M18").Value = 0.005235602094240838
$ws.Range("O18").Value = 0.09947643979057591
$ws.Range("S18").Value = 0.04712041884816754
$ws.Range("F19").Value = 0.01649646504320503
$ws.Range("H19").Value = 0.1940298507462687
$ws.Range("I19").Value = 0.09347996857816182
$ws.Range("J19").Value = 0.391987431264729
$ws.Range("K19").Value = 0.1115475255302435
$ws.Range("M19").Value = 0.0180675569520817
$ws.Range("N19").Value = 0.001571091908876669
$ws.Range("O19").Value = 0.09190887666928516
$ws.Range("S19").Value = 0.08091123330714847
